# multi lib protocol added
# Clear out the example quantity figures on the "End point" sheet (rows 17-23)
# while keeping their existing number formatting, and leave the active
# selection on the L17:L23 range as it was left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("End point")

# Remove the old sample values (20000,10000,5000,4000,3000,2000,1000 ...)
# from columns B:D, F:H and J:L for rows 17 through 23. Column E/I (blank
# spacer cells) and column A (labels) / M (date x-values) are untouched.
$ws.Range("B17:D23").ClearContents()
$ws.Range("F17:H23").ClearContents()
$ws.Range("J17:L23").ClearContents()

# Match the saved cursor/selection location recorded in the sheet view.
[void]$ws.Range("L17:L23").Select()
